# Update crypto price/volume snapshot (D:price, E:1h volume %) for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.985.58"
$ws.Range("E2").Value = "  -4.11%  "
$ws.Range("D3").Value = "3.300.47"
$ws.Range("E3").Value = "  -4.17%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.00"
$ws.Range("E5").Value = "  -3.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.39"
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.483"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.407"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "3.881.12"
$ws.Range("E12").Value = "  -3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.129"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.26"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "3.340.06"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000165"
$ws.Range("E16").Value = "  -3.42%  "
$ws.Range("D17").Value = "60.113.97"
$ws.Range("E17").Value = "  -3.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.17"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.24"
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.66"
$ws.Range("E20").Value = "  -3.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.40"
$ws.Range("E21").Value = "  -3.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.32"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.549"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "3.484.45"
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000105"
$ws.Range("E26").Value = "  -8.24%  "
$ws.Range("E27").Value = "  -5.96%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  -5.48%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.05"
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.58"
$ws.Range("E32").Value = "  -5.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.54"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.18"
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.73"
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("E37").Value = "  -6.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.72"
$ws.Range("E38").Value = "  -3.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.48"
$ws.Range("E39").Value = "  -13.88%  "
$ws.Range("D40").Value = "3.345.20"
$ws.Range("E40").Value = "  -3.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0733"
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.89"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.750"
$ws.Range("E43").Value = "  -4.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.20"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("E46").Value = "  -5.76%  "
$ws.Range("D47").Value = "2.403.93"
$ws.Range("E47").Value = "  -6.58%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.57"
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.73"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0257"
$ws.Range("E51").Value = "  -3.22%  "
